$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the header style used by the other header cells (B1:H1), then set text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I2:J75 values (row index corresponds to worksheet row)
$values = @(
    @(5, 6),
    @(6, 6),
    @(5, 5),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(5, 5),
    @(8, 8),
    @(10, 10),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(5, 5),
    @(8, 8),
    @(5, 5),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(4, 4),
    @(7, 7),
    @(5, 5),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(5, 6),
    @(6, 6),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(5, 5),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(9, 9),
    @(5, 5),
    @(4, 4),
    @(7, 7),
    @(8, 8),
    @(8, 9),
    @(3, 3),
    @(8, 8),
    @(5, 5),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(6, 7),
    @(5, 5)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
